# ============================================================
# Sync automatico del tracker - 2025-09-16 18:46:18 UTC
# Marks rows 141-144 as Completed with results, and appends
# newly scraped fixtures for 2025-09-17 as rows 150-159.
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Resolve pending predictions for 2025-09-15 / 2025-09-16 fixtures ---
$ws.Range("L141").Value = "Completed"
$ws.Range("M141").Value = "Draw"
$ws.Range("N141").Value = "Fallo"
$ws.Range("O141").Value = -1
$ws.Range("P141").Value = -100
$ws.Range("Q141").Value = "2025-09-16 04:26:19"

$ws.Range("L142").Value = "Completed"
$ws.Range("M142").Value = "Home Win"
$ws.Range("N142").Value = "Acierto"
$ws.Range("O142").Value = 1.61
$ws.Range("P142").Value = 62
$ws.Range("Q142").Value = "2025-09-16 04:26:19"

$ws.Range("L143").Value = "Completed"
$ws.Range("M143").Value = "Draw"
$ws.Range("N143").Value = "Fallo"
$ws.Range("O143").Value = -2.6
$ws.Range("P143").Value = -100
$ws.Range("Q143").Value = "2025-09-16 04:26:19"

$ws.Range("L144").Value = "Completed"
$ws.Range("M144").Value = "Home Win"
$ws.Range("N144").Value = "Acierto"
$ws.Range("O144").Value = 1.08
$ws.Range("P144").Value = 120
$ws.Range("Q144").Value = "2025-09-16 04:26:19"

# --- 2) Append newly scraped fixtures (2025-09-17) ---
# Columns A, F and H hold text that LOOKS numeric ("yyyy-MM-dd" dates
# and "NN.NN%" probabilities) but must stay plain text, exactly like
# every other row in this sheet -- format the cells as Text first so
# Excel does not silently coerce them into date serials / percentages.
$ws.Range("A150:A159").NumberFormat = "@"
$ws.Range("F150:F159").NumberFormat = "@"
$ws.Range("H150:H159").NumberFormat = "@"

# Row 150: Olympiakos Piraeus vs Pafos (UEFA Champions League)
$ws.Range("A150").Value = "2025-09-17"
$ws.Range("B150").Value = "UEFA Champions League"
$ws.Range("C150").Value = "Olympiakos Piraeus"
$ws.Range("D150").Value = "Pafos"
$ws.Range("E150").Value = "Home Win"
$ws.Range("F150").Value = "90.20%"
$ws.Range("G150").Value = 1.4
$ws.Range("H150").Value = "25.02%"
$ws.Range("I150").Value = 2.9
$ws.Range("J150").Value = 0.05
$ws.Range("K150").Value = 0.6570343192575333
$ws.Range("L150").Value = "Pending"

# Row 151: Slavia Praha vs Bodo/Glimt (UEFA Champions League)
$ws.Range("A151").Value = "2025-09-17"
$ws.Range("B151").Value = "UEFA Champions League"
$ws.Range("C151").Value = "Slavia Praha"
$ws.Range("D151").Value = "Bodo/Glimt"
$ws.Range("E151").Value = "Home Win"
$ws.Range("F151").Value = "73.18%"
$ws.Range("G151").Value = 1.75
$ws.Range("H151").Value = "26.78%"
$ws.Range("I151").Value = 2.2
$ws.Range("J151").Value = 0.03742001338783531
$ws.Range("K151").Value = 0.3742001338783531
$ws.Range("L151").Value = "Pending"

# Row 152: Fatih Karagümrük vs Istanbul Basaksehir (Süper Lig)
$ws.Range("A152").Value = "2025-09-17"
$ws.Range("B152").Value = "Süper Lig"
$ws.Range("C152").Value = "Fatih Karagümrük"
$ws.Range("D152").Value = "Istanbul Basaksehir"
$ws.Range("E152").Value = "Away Win"
$ws.Range("F152").Value = "64.52%"
$ws.Range("G152").Value = 1.91
$ws.Range("H152").Value = "22.00%"
$ws.Range("I152").Value = 1.5
$ws.Range("J152").Value = 0.0255302801976002
$ws.Range("K152").Value = 0.255302801976002
$ws.Range("L152").Value = "Pending"

# Row 153: Samsunspor vs Kasimpasa (Süper Lig)
$ws.Range("A153").Value = "2025-09-17"
$ws.Range("B153").Value = "Süper Lig"
$ws.Range("C153").Value = "Samsunspor"
$ws.Range("D153").Value = "Kasimpasa"
$ws.Range("E153").Value = "Home Win"
$ws.Range("F153").Value = "70.54%"
$ws.Range("G153").Value = 1.8
$ws.Range("H153").Value = "25.70%"
$ws.Range("I153").Value = 1.9
$ws.Range("J153").Value = 0.03371245248403483
$ws.Range("K153").Value = 0.3371245248403483
$ws.Range("L153").Value = "Pending"

# Row 154: Genk vs Charleroi (Jupiler Pro League)
$ws.Range("A154").Value = "2025-09-17"
$ws.Range("B154").Value = "Jupiler Pro League"
$ws.Range("C154").Value = "Genk"
$ws.Range("D154").Value = "Charleroi"
$ws.Range("E154").Value = "Home Win"
$ws.Range("F154").Value = "72.78%"
$ws.Range("G154").Value = 1.7
$ws.Range("H154").Value = "22.50%"
$ws.Range("I154").Value = 2
$ws.Range("J154").Value = 0.03390349092324032
$ws.Range("K154").Value = 0.3390349092324031
$ws.Range("L154").Value = "Pending"

# Row 155: Ajax vs Inter (UEFA Champions League)
$ws.Range("A155").Value = "2025-09-17"
$ws.Range("B155").Value = "UEFA Champions League"
$ws.Range("C155").Value = "Ajax"
$ws.Range("D155").Value = "Inter"
$ws.Range("E155").Value = "Away Win"
$ws.Range("F155").Value = "71.12%"
$ws.Range("G155").Value = 1.8
$ws.Range("H155").Value = "26.74%"
$ws.Range("I155").Value = 2
$ws.Range("J155").Value = 0.03501982061867995
$ws.Range("K155").Value = 0.3501982061867995
$ws.Range("L155").Value = "Pending"

# Row 156: Paris Saint Germain vs Atalanta (UEFA Champions League)
$ws.Range("A156").Value = "2025-09-17"
$ws.Range("B156").Value = "UEFA Champions League"
$ws.Range("C156").Value = "Paris Saint Germain"
$ws.Range("D156").Value = "Atalanta"
$ws.Range("E156").Value = "Home Win"
$ws.Range("F156").Value = "87.75%"
$ws.Range("G156").Value = 1.45
$ws.Range("H156").Value = "25.97%"
$ws.Range("I156").Value = 2.9
$ws.Range("J156").Value = 0.05
$ws.Range("K156").Value = 0.6053342623591915
$ws.Range("L156").Value = "Pending"

# Row 157: Liverpool vs Atletico Madrid (UEFA Champions League)
$ws.Range("A157").Value = "2025-09-17"
$ws.Range("B157").Value = "UEFA Champions League"
$ws.Range("C157").Value = "Liverpool"
$ws.Range("D157").Value = "Atletico Madrid"
$ws.Range("E157").Value = "Home Win"
$ws.Range("F157").Value = "84.57%"
$ws.Range("G157").Value = 1.5
$ws.Range("H157").Value = "25.58%"
$ws.Range("I157").Value = 2.9
$ws.Range("J157").Value = 0.05
$ws.Range("K157").Value = 0.5369698608511987
$ws.Range("L157").Value = "Pending"

# Row 158: Bayern München vs Chelsea (UEFA Champions League)
$ws.Range("A158").Value = "2025-09-17"
$ws.Range("B158").Value = "UEFA Champions League"
$ws.Range("C158").Value = "Bayern München"
$ws.Range("D158").Value = "Chelsea"
$ws.Range("E158").Value = "Home Win"
$ws.Range("F158").Value = "78.13%"
$ws.Range("G158").Value = 1.65
$ws.Range("H158").Value = "27.62%"
$ws.Range("I158").Value = 2.6
$ws.Range("J158").Value = 0.04447810602235108
$ws.Range("K158").Value = 0.4447810602235108
$ws.Range("L158").Value = "Pending"

# Row 159: New York City FC vs Columbus Crew (Major League Soccer)
$ws.Range("A159").Value = "2025-09-17"
$ws.Range("B159").Value = "Major League Soccer"
$ws.Range("C159").Value = "New York City FC"
$ws.Range("D159").Value = "Columbus Crew"
$ws.Range("E159").Value = "Home Win"
$ws.Range("F159").Value = "55.20%"
$ws.Range("G159").Value = 2.1
$ws.Range("H159").Value = "14.76%"
$ws.Range("I159").Value = 0.8
$ws.Range("J159").Value = 0.01446775813821211
$ws.Range("K159").Value = 0.1446775813821211
$ws.Range("L159").Value = "Pending"

# Drop the Text-format override again so the new rows keep the same
# unstyled ("Normal") look as every pre-existing data row.
$ws.Range("A150:Q159").Style = "Normal"

